$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.456.89'
$ws.Range('E2').Value = '  +1.48%  '
$ws.Range('D3').Value = '1.676.71'
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5308'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.31%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2694'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06414'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.59%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07827'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.22%  '
$ws.Range('D12').Value = '1.684.29'
$ws.Range('E12').Value = '  +3.04%  '
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5580'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.68%  '
$ws.Range('D15').Value = '0.0₅8353'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.75'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.19%  '
$ws.Range('D17').Value = '26.495.64'
$ws.Range('E17').Value = '  +1.65%  '
$ws.Range('E18').Value = '  -0.10%  '
$ws.Range('E19').Value = '  +1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.11'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.06%  '
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.347'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.24%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '142.63'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1288'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.396'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.441'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.06312'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.273'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.58%  '
$ws.Range('E31').Value = '  +5.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.455'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.55%  '
$ws.Range('E33').Value = '  +2.57%  '
$ws.Range('E34').Value = '  +2.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6201'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +8.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.427'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.791'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.165'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.46%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01636'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.24%  '
$ws.Range('D40').Value = '1.087.09'
$ws.Range('E40').Value = '  +4.97%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8660'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.53%  '
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.41'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('D44').Value = '1.822.03'
$ws.Range('E44').Value = '  +1.98%  '
$ws.Range('E45').Value = '  +3.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.194'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.72%  '
$ws.Range('E47').Value = '  +0.01%  '
$ws.Range('D48').Value = '0.0₈103'
$ws.Range('E48').Value = '  -3.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05208'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.77%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.483'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.95%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.040'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.12%  '
